$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header strings: "_old" -> "_FV2404", "_new" -> "_FV2410" ---
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_old$", "_FV2404")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_new$", "_FV2410")
}

# --- Turn the data range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U74")
$lo = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# --- Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)

Write-Host "Edit complete"
